# staging smoke test debugging - MSRP price corrections for rows 29-34
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29 (9700): MSRP 53000 -> 53100
$ws.Range("D29").Value = 53100

# Row 30 (9700PM): MSRP 55790 -> 55890
$ws.Range("D30").Value = 55890

# Row 31 (9710): MSRP 64265 -> 64365
$ws.Range("D31").Value = 64365

# Row 32 (9625): MSRP 86480 -> 86580, freight 1295 -> 1025
$ws.Range("D32").Value = 86580
$ws.Range("E32").Value = 1025

# Row 33 (9620): MSRP 91480 -> 91580, freight 1295 -> 1025
$ws.Range("D33").Value = 91580
$ws.Range("E33").Value = 1025

# Row 34 (9620 (SE)): was blank placeholder, now a real MSRP of 99310, freight 1295 -> 1025
$ws.Range("D34").Value = 99310
$ws.Range("D34").NumberFormat = "#,##0"
$ws.Range("E34").Value = 1025

# leave selection where the edits ended up
$ws.Range("D29").Select()
